$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1246.91
$ws.Range("I132").Value = 1163.7609
$ws.Range("J132").Value = 2203.125
$ws.Range("K132").Value = 3491.2827
$ws.Range("L132").Value = 6609.375
$ws.Range("M132").Value = -961.2826999999997
$ws.Range("N132").Value = -11669.375
$ws.Range("H137").Value = 4442.5
$ws.Range("I137").Value = 4442.5
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 13327.5
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -10777.5
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 163383.77
$ws.Range("I138").Value = 33210.484
$ws.Range("J138").Value = 272447.88
$ws.Range("K138").Value = 99631.45199999999
$ws.Range("L138").Value = 817343.64
$ws.Range("M138").Value = -94491.45199999999
$ws.Range("N138").Value = -827623.64

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3043.2698
$ws.Range("I32").Value = 2909.885
$ws.Range("J32").Value = 8845.5
$ws.Range("K32").Value = 2909.885
$ws.Range("L32").Value = 8845.5
$ws.Range("M32").Value = -2622.885
$ws.Range("N32").Value = -9419.5
$ws.Range("H45").Value = 1465.7097
$ws.Range("I45").Value = 1289.9231
$ws.Range("J45").Value = 2379.8
$ws.Range("K45").Value = 1289.9231
$ws.Range("L45").Value = 2379.8
$ws.Range("M45").Value = -912.9231
$ws.Range("N45").Value = -3133.8
$ws.Range("H52").Value = 99999
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 99999
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 99999
$ws.Range("N52").Value = -100635
$ws.Range("H98").Value = 49660
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 49660
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 49660
$ws.Range("N98").Value = -55650
$ws.Range("H101").Value = 80326.664
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 80326.664
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 80326.664
$ws.Range("N101").Value = -86816.664
$ws.Range("H106").Value = 85934
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 85934
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 85934
$ws.Range("N106").Value = -88458
$ws.Range("H120").Value = 89709.5
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 89709.5
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 89709.5
$ws.Range("N120").Value = -99385.5
$ws.Range("H121").Value = 85994.8
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 85994.8
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 85994.8
$ws.Range("N121").Value = -89488.8
$ws.Range("H122").Value = 1515.1428
$ws.Range("I122").Value = 1335.5385
$ws.Range("J122").Value = 3850
$ws.Range("K122").Value = 4006.6155
$ws.Range("L122").Value = 11550
$ws.Range("M122").Value = -1556.6155
$ws.Range("N122").Value = -16450
$ws.Range("H132").Value = 4832.091
$ws.Range("I132").Value = 2981.0625
$ws.Range("J132").Value = 9768.166999999999
$ws.Range("K132").Value = 8943.1875
$ws.Range("L132").Value = 29304.501
$ws.Range("M132").Value = -6413.1875
$ws.Range("N132").Value = -34364.501

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 76570.664
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 76570.664
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 76570.664
$ws.Range("N6").Value = -76796.664
$ws.Range("H13").Value = 98831.664
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 98831.664
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 98831.664
$ws.Range("N13").Value = -99167.664
$ws.Range("H29").Value = 11707.75
$ws.Range("I29").Value = 10610.333
$ws.Range("J29").Value = 15000
$ws.Range("K29").Value = 10610.333
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = -10321.333
$ws.Range("N29").Value = -15578
$ws.Range("H36").Value = 1641.2
$ws.Range("I36").Value = 1641.2
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1641.2
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1107.2
$ws.Range("H51").Value = 99662
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 99662
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 99662
$ws.Range("N51").Value = -100644
$ws.Range("H52").Value = 77979
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 77979
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 77979
$ws.Range("N52").Value = -78505
$ws.Range("H119").Value = 72705.5
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 72705.5
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 72705.5
$ws.Range("N119").Value = -82381.5
$ws.Range("H120").Value = 122170.336
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 122170.336
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 122170.336
$ws.Range("N120").Value = -131846.336
$ws.Range("H121").Value = 77979
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 77979
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 77979
$ws.Range("N121").Value = -81473
$ws.Range("H134").Value = 3196.12
$ws.Range("I134").Value = 1545.0883
$ws.Range("J134").Value = 4565.268
$ws.Range("K134").Value = 4635.2649
$ws.Range("L134").Value = 13695.804
$ws.Range("M134").Value = -2100.2649
$ws.Range("N134").Value = -18765.804

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 20000
$ws.Range("I41").Value = 20000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 20000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -19572
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H51").Value = 94350
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 94350
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 94350
$ws.Range("N51").Value = -95822
$ws.Range("H61").Value = 94350
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 94350
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 94350
$ws.Range("N61").Value = -95046
$ws.Range("H107").Value = 1547.122
$ws.Range("I107").Value = 983.5333000000001
$ws.Range("J107").Value = 3084.182
$ws.Range("K107").Value = 983.5333000000001
$ws.Range("L107").Value = 3084.182
$ws.Range("M107").Value = 936.4666999999999
$ws.Range("N107").Value = -6924.182

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4106.7
$ws.Range("I81").Value = 1975.75
$ws.Range("J81").Value = 5527.3335
$ws.Range("K81").Value = 5927.25
$ws.Range("L81").Value = 16582.0005
$ws.Range("M81").Value = -4804.25
$ws.Range("N81").Value = -18828.0005
$ws.Range("H84").Value = 4106.7
$ws.Range("I84").Value = 1975.75
$ws.Range("J84").Value = 5527.3335
$ws.Range("K84").Value = 17781.75
$ws.Range("L84").Value = 49746.0015
$ws.Range("M84").Value = -12165.75
$ws.Range("N84").Value = -60978.0015

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 108.57143
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 167.5
$ws.Range("K2").Value = 30
$ws.Range("L2").Value = 167.5
$ws.Range("M2").Value = 83
$ws.Range("N2").Value = -393.5
$ws.Range("H80").Value = 1800.68
$ws.Range("I80").Value = 1815.9333
$ws.Range("J80").Value = 1777.8
$ws.Range("K80").Value = 1815.9333
$ws.Range("L80").Value = 1777.8
$ws.Range("M80").Value = -817.9332999999999
$ws.Range("N80").Value = -3773.8
$ws.Range("H83").Value = 1800.68
$ws.Range("I83").Value = 1815.9333
$ws.Range("J83").Value = 1777.8
$ws.Range("K83").Value = 9079.666499999999
$ws.Range("L83").Value = 8889
$ws.Range("M83").Value = -4087.666499999999
$ws.Range("N83").Value = -18873
$ws.Range("H97").Value = 924.65625
$ws.Range("I97").Value = 817.53845
$ws.Range("J97").Value = 997.9474
$ws.Range("K97").Value = 817.53845
$ws.Range("L97").Value = 997.9474
$ws.Range("M97").Value = -321.53845
$ws.Range("N97").Value = -1989.9474
$ws.Range("H102").Value = 29536.166
$ws.Range("I102").Value = 1861.9412
$ws.Range("J102").Value = 499998
$ws.Range("K102").Value = 1861.9412
$ws.Range("L102").Value = 499998
$ws.Range("M102").Value = -239.9412
$ws.Range("N102").Value = -503242
$ws.Range("H107").Value = 334.45456
$ws.Range("I107").Value = 238.06667
$ws.Range("J107").Value = 541
$ws.Range("K107").Value = 238.06667
$ws.Range("L107").Value = 541
$ws.Range("M107").Value = 1681.93333
$ws.Range("N107").Value = -4381

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 130.90909
$ws.Range("I55").Value = 119
$ws.Range("J55").Value = 250
$ws.Range("K55").Value = 119
$ws.Range("L55").Value = 250
$ws.Range("M55").Value = 54
$ws.Range("N55").Value = -596
$ws.Range("H59").Value = 63441.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 63441.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 63441.5
$ws.Range("N59").Value = -64749.5
$ws.Range("H110").Value = 80000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 80000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180
$ws.Range("H112").Value = 74582
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 74582
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 74582
$ws.Range("N112").Value = -77536
$ws.Range("H117").Value = 91841.336
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 91841.336
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 91841.336
$ws.Range("N117").Value = -101019.336
$ws.Range("H120").Value = 79999
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 79999
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 79999
$ws.Range("N120").Value = -89675
$ws.Range("H121").Value = 63332.668
$ws.Range("I121").Value = 99999
$ws.Range("J121").Value = 44999.5
$ws.Range("K121").Value = 99999
$ws.Range("L121").Value = 44999.5
$ws.Range("M121").Value = -98252
$ws.Range("N121").Value = -48493.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H114").Value = 49999.855
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 49999.855
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 49999.855
$ws.Range("N114").Value = -58677.855
$ws.Range("H122").Value = 2357.8865
$ws.Range("I122").Value = 2150.2432
$ws.Range("J122").Value = 3455.4285
$ws.Range("K122").Value = 6450.7296
$ws.Range("L122").Value = 10366.2855
$ws.Range("M122").Value = -4000.7296
$ws.Range("N122").Value = -15266.2855
$ws.Range("H132").Value = 1828.6346
$ws.Range("I132").Value = 1930.4186
$ws.Range("J132").Value = 1342.3334
$ws.Range("K132").Value = 5791.2558
$ws.Range("L132").Value = 4027.0002
$ws.Range("M132").Value = -3261.2558
$ws.Range("N132").Value = -9087.0002
